$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Keyword in row 5 from "geaca" to "blugi"
$ws.Range("A5").Value = "blugi"

# Update the selected/active cell to A5
$ws.Range("A5").Select()
